# Update the "Training Dashboard" sheet with the latest progress snapshot
# dated 04-Nov-2025 (previously 03-Nov-2025).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# The "LAST UPDATE" column stores plain text dates (not real date values).
# A bare Value assignment of a date-looking string makes Excel auto-convert
# it to a real date serial (and change the cell's number format), so write
# it with a leading apostrophe to force literal text, then restore the
# cell's original look by re-pasting the (unchanged) formatting from its
# same-styled neighbour in column H.

# Row 3: Endangered by Electricity A safety Training (SOPs)
$ws.Range("H3").Value = -100
$ws.Range("I3").Value = "'04-Nov-2025"
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)

# Row 4: Incident Escalation Process (SOPs)
$ws.Range("H4").Value = 251
$ws.Range("I4").Value = "'04-Nov-2025"
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)

# Row 5: ISO 55001 (Other Trainings)
$ws.Range("H5").Value = 286
$ws.Range("I5").Value = "'04-Nov-2025"
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)

$excel.CutCopyMode = $false
